$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.491.85"
$ws.Range("E2").Value = "  -3.45%  "

$ws.Range("D3").Value = "2.480.01"
$ws.Range("E3").Value = "  -5.98%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.53%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  -3.29%  "

$ws.Range("D9").Value = "2.478.25"
$ws.Range("E9").Value = "  -5.96%  "

$ws.Range("E10").Value = "  -8.39%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.87%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.154"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.41%  "

$ws.Range("E13").Value = "  -6.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.61%  "

$ws.Range("D15").Value = "2.927.10"
$ws.Range("E15").Value = "  -6.04%  "

$ws.Range("E16").Value = "  -8.54%  "

$ws.Range("D17").Value = "61.370.41"
$ws.Range("E17").Value = "  -3.59%  "

$ws.Range("D18").Value = "2.482.33"
$ws.Range("E18").Value = "  -6.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.66%  "

$ws.Range("E21").Value = "  -7.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "322.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.45%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.68%  "

$ws.Range("D26").Value = "0.0₃0998"
$ws.Range("E26").Value = "  -8.33%  "

$ws.Range("D27").Value = "2.608.47"
$ws.Range("E27").Value = "  -5.73%  "

$ws.Range("E28").Value = "  -5.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "545.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.38%  "

$ws.Range("E32").Value = "  -3.80%  "

$ws.Range("E33").Value = "  -5.33%  "

$ws.Range("E34").Value = "  -7.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.97%  "

$ws.Range("E36").Value = "  -9.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.04%  "

$ws.Range("E39").Value = "  -4.77%  "

$ws.Range("E40").Value = "  -5.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "147.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.41%  "

$ws.Range("E42").Value = "  -7.80%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.53%  "

$ws.Range("E45").Value = "  -6.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.82%  "

$ws.Range("E49").Value = "  -7.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.598"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.48%  "

$ws.Range("E51").Value = "  -4.99%  "
